# matlab week 2: almost done with RLE decoder
# Shrink the Game-of-Life input grid from 9x8 (A1:I8) down to 5x5 (A1:E5)
# and set it to a simple vertical 3-cell "blinker" pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything first so the old 9x8 footprint (columns F:I, rows 6:8)
# doesn't linger outside the new A1:E5 used range.
$ws.Cells.Clear()

# New 5x5 grid values (row-major, A..E across, 1..5 down).
$grid = @(
    @(0, 0, 0, 0, 0),
    @(0, 0, 1, 0, 0),
    @(0, 0, 1, 0, 0),
    @(0, 0, 1, 0, 0),
    @(0, 0, 0, 0, 0)
)

for ($r = 1; $r -le 5; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value = $grid[$r - 1][$c - 1]
    }
}

# Match the saved selection from the edit (cell D3).
$ws.Range("D3").Select()
